$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5854969999999999
$ws.Range("H2").Value = 1.170994
$ws.Range("M2").Value = 30.801072
$ws.Range("N2").Value = 61.602144
$ws.Range("O2").Value = 0.5373480691764108
$ws.Range("P2").Value = 0.485871843331092
$ws.Range("Q2").Value = 18.033935252784
$ws.Range("R2").Value = 72.13574101113599
$ws.Range("S2").Value = 0.5373480691764108
$ws.Range("T2").Value = 0.485871843331092

# Update row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5854969999999999
$ws.Range("H3").Value = 1.170994
$ws.Range("O3").Value = 0.2096145064786482
$ws.Range("P3").Value = 0.2843011610923331
$ws.Range("Q3").Value = 7.034871165860999
$ws.Range("R3").Value = 42.209226995166
$ws.Range("S3").Value = 0.2096145064786482
$ws.Range("T3").Value = 0.2843011610923331

# Update row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5854969999999999
$ws.Range("H4").Value = 1.170994
$ws.Range("M4").Value = 0.1305583333333333
$ws.Range("N4").Value = 0.391675
$ws.Range("O4").Value = 0.002277689176907768
$ws.Range("P4").Value = 0.003089240761436898
$ws.Range("Q4").Value = 0.07644151249166665
$ws.Range("R4").Value = 0.45864907495
$ws.Range("S4").Value = 0.002277689176907768
$ws.Range("T4").Value = 0.003089240761436898

# Update row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5854969999999999
$ws.Range("H5").Value = 1.170994
$ws.Range("M5").Value = 14.373679
$ws.Range("N5").Value = 28.747358
$ws.Range("O5").Value = 0.2507597351680332
$ws.Range("P5").Value = 0.2267377548151379
$ws.Range("Q5").Value = 8.415745933463
$ws.Range("R5").Value = 33.662983733852
$ws.Range("S5").Value = 0.2507597351680332
$ws.Range("T5").Value = 0.2267377548151379

# Delete rows 6 and 7 (Neutrophils / Resolving-Mac target-cluster rows)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
